# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q3" (i.e. right
#    before "总计"), populate it with the new quarter's fund-holding rows,
#    matching the look (bold+bordered header / first column) of the other
#    quarter sheets.
# 2. Insert a new summary row at the top of the "总计" sheet's data
#    (row 2, pushing the existing rows down) for the 2022-Q1 totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2021-Q3")
$q1new = $wb.Worksheets.Add($null, $q3)
$q1new.Name = "2022-Q1"

# Header row (bold font + thin border all round, centered/top aligned -
# same look as the other quarter sheets)
$header = $q1new.Range("B1:H1")
$header.Value = "" # no-op placeholder so style application below has a range to work on
$q1new.Range("B1").Value = "基金代码"
$q1new.Range("C1").Value = "基金名称"
$q1new.Range("D1").Value = "基金规模"
$q1new.Range("E1").Value = "股票总仓位"
$q1new.Range("F1").Value = "仓位占比"
$q1new.Range("G1").Value = "持有市值(亿元)"
$q1new.Range("H1").Value = "仓位排名"
$header.Font.Bold = $true
$header.HorizontalAlignment = "Center"
$header.VerticalAlignment = "Top"
$header.Borders.LineStyle = "Continuous"
$header.Borders.Weight = "Thin"

# Data rows - column A is the numeric row index (0-based), columns D:G
# are free-text numbers (percentages / amounts kept as text, not
# converted), column H is numeric.
$rows = @(
    @{ idx=0; code="512200"; name="南方中证全指房地产ETF";     scale="28.63"; pos="99.85"; pct="1.69"; value="0.4838"; rank=10 },
    @{ idx=1; code="160218"; name="国泰国证房地产行业指数";     scale="6.35";  pos="95.04"; pct="1.85"; value="0.1175"; rank=10 },
    @{ idx=2; code="160628"; name="鹏华中证800地产指数（LOF）"; scale="3.38";  pos="94.35"; pct="2.39"; value="0.0808"; rank=10 },
    @{ idx=3; code="515060"; name="华夏中证全指房地产ETF";     scale="2.37";  pos="98.82"; pct="1.66"; value="0.0393"; rank=10 }
)

$r = 2
foreach ($row in $rows) {
    $aCell = $q1new.Range("A$r")
    $aCell.Value = $row.idx
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = "Center"
    $aCell.VerticalAlignment = "Top"
    $aCell.Borders.LineStyle = "Continuous"
    $aCell.Borders.Weight = "Thin"

    $q1new.Range("B$r").Value = $row.code
    $q1new.Range("C$r").Value = $row.name

    $textRange = $q1new.Range("D$r:G$r")
    $textRange.NumberFormat = "@"
    $q1new.Range("D$r").Value = $row.scale
    $q1new.Range("E$r").Value = $row.pos
    $q1new.Range("F$r").Value = $row.pct
    $q1new.Range("G$r").Value = $row.value

    $q1new.Range("H$r").Value = $row.rank

    $r = $r + 1
}

$q1new.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Update "总计" sheet: insert a new 2022-Q1 summary row above the
#    existing rows (so the list stays newest-quarter-first).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$newRowA = $totalSheet.Range("A2")
$newRowA.Value = 0
$newRowA.Font.Bold = $true
$newRowA.HorizontalAlignment = "Center"
$newRowA.VerticalAlignment = "Top"
$newRowA.Borders.LineStyle = "Continuous"
$newRowA.Borders.Weight = "Thin"

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.72

# Renumber the index column (A) for the rows that got pushed down so it
# keeps reading 0,1,2,3 top to bottom.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
